# Fruta / hortaliza, semanal
# Insert a new weekly observation row after the current row 244, shifting
# every subsequent row down by one (the existing rows 245-260 become 246-261).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 245 (row 244 stays put, old 245.. shift down).
$ws.Rows.Item(245).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Range("A245").Value = 1
$ws.Range("B245").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C245").Value = "Arica y Parinacota"
$ws.Range("D245").Value = 44746
$ws.Range("E245").Value = 15
$ws.Range("F245").Value = "Fruta"
$ws.Range("G245").Value = 100108
$ws.Range("H245").Value = "Tropicales y subtropicales"
$ws.Range("I245").Value = 100108006
$ws.Range("J245").Value = "Plátano"
$ws.Range("K245").Value = "Sin especificar"
$ws.Range("L245").Value = "Pintón"
$ws.Range("M245").Value = 120
$ws.Range("N245").Value = 19000
$ws.Range("O245").Value = 20000
$ws.Range("P245").Value = 19500
$ws.Range("Q245").Value = "$/caja 20 kilos"
$ws.Range("R245").Value = "Ecuador"
$ws.Range("S245").Value = 975
$ws.Range("T245").Value = 20
